# Refresh the cryptos table: update Price/Volume(1h) figures and, where the
# source ranking reordered two adjacent coins, swap their Coin/Link/Price/Volume.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.749.28"
$ws.Range("E2").Value = "  +4.64%  "

$ws.Range("D3").Value = "2.280.79"
$ws.Range("E3").Value = "  +2.46%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'231.95"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("D7").Value = "'61.57"
$ws.Range("E7").Value = "  +1.59%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "'0.421"
$ws.Range("E9").Value = "  +5.14%  "

$ws.Range("D10").Value = "'0.0926"
$ws.Range("E10").Value = "  +3.87%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").Value = "2.622.15"
$ws.Range("E12").Value = "  +2.63%  "

$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("D14").Value = "'23.71"
$ws.Range("E14").Value = "  +9.18%  "

$ws.Range("D15").Value = "'5.75"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").Value = "'0.810"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("D17").Value = "2.295.61"
$ws.Range("E17").Value = "  +3.37%  "

$ws.Range("D18").Value = "43.547.49"
$ws.Range("E18").Value = "  +4.49%  "

$ws.Range("D19").Value = "0.0₃0935"
$ws.Range("E19").Value = "  +4.88%  "

$ws.Range("D20").Value = "'73.27"
$ws.Range("E20").Value = "  +0.73%  "

$ws.Range("E21").Value = "  +3.65%  "

$ws.Range("D22").Value = "'250.87"
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("E24").Value = "  +7.34%  "

$ws.Range("D25").Value = "'2.37"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").Value = "'9.86"
$ws.Range("E26").Value = "  +3.12%  "

$ws.Range("D27").Value = "'170.40"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("D29").Value = "'20.62"
$ws.Range("E29").Value = "  +3.47%  "

$ws.Range("D30").Value = "'1.49"
$ws.Range("E30").Value = "  +6.56%  "

$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'5.05"
$ws.Range("E33").Value = "  +2.16%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.77"
$ws.Range("E34").Value = "  +3.49%  "

$ws.Range("D35").Value = "'0.0661"
$ws.Range("E35").Value = "  +5.89%  "

$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "'6.54"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.43"
$ws.Range("E37").Value = "  +3.28%  "

$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("D39").Value = "'0.0252"
$ws.Range("E39").Value = "  +4.70%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "'8.88"
$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "'0.000219"
$ws.Range("E42").Value = "  -14.44%  "

$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "'4.53"
$ws.Range("E43").Value = "  -5.51%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'1.22"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0970"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("D46").Value = "'98.44"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").Value = "1.474.29"
$ws.Range("E47").Value = "  +0.35%  "

$ws.Range("D48").Value = "'16.70"
$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("E49").Value = "  +9.59%  "

$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "'2.78"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'1.09"
$ws.Range("E51").Value = "  +1.43%  "
